$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Hyperlinks.Delete()
Write-Host ("hyperlinks count: " + $ws.Hyperlinks.Count)
